{"js": "// Office.js (Word JavaScript API) script.\n// This is the body of: async (context) => { ... }\n//\n// Implements the \"Changed technology stack section in documentation\" edit:\n//\n//   1. In the paragraph describing the tech stack, the sentence\n//        \"...a Loopback based API backend.\"\n//      becomes\n//        \"...a Loopback based API framework for the backend.\"\n//      i.e. \" framework for the\" is inserted right after \"...API\".\n//\n//   2. In the Source-Control paragraph, \"GitHub\" + \" \" (a run boundary\n//      right after the word \"GitHub\") is normalized into a single\n//      \"GitHub \" text run (purely a run-structure clean-up; the visible\n//      text is unchanged).\n\n// --- 1. Tech stack sentence -----------------------------------------\nconst apiResults = context.document.body.search(\n  \"Loopback based API\",\n  { matchCase: true, matchWholeWord: false }\n);\napiResults.load(\"items\");\nawait context.sync();\n\nif (apiResults.items.length > 0) {\n  // Insert the new phrase immediately after the matched text (\"...API\"),\n  // right before \" backend.\" that follows it.\n  apiResults.items[0].insertText(\" framework for the\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// --- 2. \"GitHub \" run normalization ----------------------------------\nconst githubResults = context.document.body.search(\n  \"GitHub \",\n  { matchCase: true, matchWholeWord: false }\n);\ngithubResults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < githubResults.items.length; i++) {\n  const match = githubResults.items[i];\n  if (match.text === \"GitHub \") {\n    match.insertText(\"GitHub \", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word / $app is the Word.Application, $word.ActiveDocument (or $d) is the document.\n#\n# Implements the \"Changed technology stack section in documentation\" edit:\n#\n#   1. In the paragraph describing the tech stack, the sentence\n#        \"...a Loopback based API backend.\"\n#      becomes\n#        \"...a Loopback based API framework for the backend.\"\n#      i.e. \" framework for the\" is inserted right after \"...API\".\n#\n#   2. In the Source-Control paragraph, \"GitHub\" + \" \" (a run boundary\n#      right after the word \"GitHub\") is normalized into a single\n#      \"GitHub \" text run (purely a run-structure clean-up; the visible\n#      text is unchanged).\n\n$d = $word.ActiveDocument\n\n# --- 1. Tech stack sentence -------------------------------------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Forward = $true\n$rng.Find.Execute(\"Loopback based API\")\nif ($rng.Find.Found) {\n    $rng.Collapse(0)            # wdCollapseEnd\n    $rng.InsertAfter(\" framework for the\")\n}\n\n# --- 2. \"GitHub \" run normalization ------------------------------------\n# \"GitHub\" appears twice in this paragraph as a whole word: once inside\n# the opening sentence (\"GitHub was used...\") and once as its own\n# standalone run right before \"repositories\". Walk the \"GitHub \" matches\n# (including the trailing space) to find that second, standalone one.\n$probe = $d.Content\n$probe.Find.ClearFormatting()\n$probe.Find.Forward = $true\n$found = $probe.Find.Execute(\"GitHub \")\n$i = 0\n$target = $null\nwhile ($found -and $i -lt 10) {\n    $i = $i + 1\n    if ($i -eq 2) {\n        $target = $d.Range($probe.Start, $probe.End)\n    }\n    $probe.Collapse(0)\n    $found = $probe.Find.Execute(\"GitHub \")\n}\n\nif ($target -ne $null) {\n    # Re-stamping the range's FormattedText merges the two runs it spans\n    # (\"GitHub\" + \" \") into a single run, without touching the\n    # neighboring runs before/after it. A same-value assignment is a\n    # no-op in this engine, so first stamp a same-length placeholder to\n    # force the edit, then set the real (unchanged-looking) text.\n    $placeholder = $target.FormattedText\n    $placeholder.Text = \"GitHubX\"\n    $target.FormattedText = $placeholder\n\n    $final = $target.FormattedText\n    $final.Text = \"GitHub \"\n    $target.FormattedText = $final\n}\n"}
